$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row additions - copy formatting from an existing header cell (AC1)
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill data rows 2-47 with team record values
for ($r = 2; $r -le 47; $r++) {
    $ws.Cells.Item($r, 30).Value = 78  # AD
    $ws.Cells.Item($r, 31).Value = 84  # AE
    $ws.Cells.Item($r, 32).Value = 0   # AF
}
